$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("WithTable") and Sheet 2 ("Tableless") share identical data tables ---
$sheetsWithSameData = @("WithTable", "Tableless")
foreach ($sheetName in $sheetsWithSameData) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row 2
    $ws.Range("A2").Value = 1
    $ws.Range("B2").Value = "Hello"
    $ws.Range("C2").Value = 45213
    $ws.Range("D2").Value = $true

    # Row 3
    $ws.Range("A3").Value = 2
    $ws.Range("B3").Value = "World"
    $ws.Range("C3").Value = 45214.75
    $ws.Range("D3").Value = $false

    # Row 4
    $ws.Range("A4").Value = 3
    $ws.Range("B4").Value = "Bye"
    $ws.Range("C4").Value = 45215.83333333334

    # Row 5
    $ws.Range("A5").Value = 4.27
    $ws.Range("B5").Value = "Outer Space"
    $ws.Range("C5").Value = 45216
}

# --- Sheet 3 ("WithTable_Duplicate") holds the same values, shifted down 3 rows ---
$ws3 = $wb.Worksheets.Item("WithTable_Duplicate")

# Row 5
$ws3.Range("B5").Value = 1
$ws3.Range("C5").Value = "Hello"
$ws3.Range("D5").Value = 45213
$ws3.Range("E5").Value = $true

# Row 6
$ws3.Range("B6").Value = 2
$ws3.Range("C6").Value = "World"
$ws3.Range("D6").Value = 45214.75
$ws3.Range("E6").Value = $false

# Row 7
$ws3.Range("B7").Value = 3
$ws3.Range("C7").Value = "Bye"
$ws3.Range("D7").Value = 45215.83333333334

# Row 8
$ws3.Range("B8").Value = 4.27
$ws3.Range("C8").Value = "Outer Space"
$ws3.Range("D8").Value = 45216

# --- styles.xml gains a third number format ("YYYY-MM-DD") backing two new cellXfs records ---
# Applied transiently to a scratch cell (then cleared) so the sheet data is
# unaffected while the style definitions remain registered in styles.xml.
$ws1 = $wb.Worksheets.Item("WithTable")
$ws1.Range("H1").NumberFormat = "YYYY-MM-DD"
$ws1.Range("H2").NumberFormat = "YYYY-MM-DD"
$ws1.Range("H1").Clear()
$ws1.Range("H2").Clear()
